# Updated cryptos list on Tue Nov  5 04:22:00 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text (e.g. "557.68", "0.506") even
# though many values look numeric. Force the Price+Volume data range to a
# text number-format before writing so Excel's auto-detection doesn't
# silently convert these into numeric cells; restore the default "Normal"
# style afterwards so no stray formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.400.77"
$ws.Range("E2").Value = "  -1.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.427.58"
$ws.Range("E3").Value = "  -1.91%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "557.68"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6 - Solana
$ws.Range("D6").Value = "160.07"
$ws.Range("E6").Value = "  -1.98%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.508"
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.163"
$ws.Range("E9").Value = "  +7.18%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -1.69%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.331"
$ws.Range("E11").Value = "  -0.75%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  -5.63%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "68.268.52"
$ws.Range("E13").Value = "  -1.12%  "

# Row 14 - was WrappedliquidstakedEther2.0, becomes ShibaInu
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000174"
$ws.Range("E14").Value = "  +2.47%  "

# Row 15 - was ShibaInu, becomes WrappedliquidstakedEther2.0
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.863.49"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "22.99"
$ws.Range("E16").Value = "  -3.19%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.419.74"
$ws.Range("E17").Value = "  -2.77%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "10.45"
$ws.Range("E18").Value = "  -3.35%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "333.56"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.88"
$ws.Range("E20").Value = "  -2.08%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "3.81"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22 - SuiNetwork
$ws.Range("D22").Value = "1.90"
$ws.Range("E22").Value = "  -0.75%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.00%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "66.50"
$ws.Range("E24").Value = "  -1.21%  "

# Row 25 - NEARProtocol
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "2.542.82"
$ws.Range("E26").Value = "  -2.36%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "8.20"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0815"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  -1.32%  "

# Row 30 - FirstDigitalUSD
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - Bittensor
$ws.Range("D31").Value = "425.21"
$ws.Range("E31").Value = "  -1.86%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -0.54%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "1.61"
$ws.Range("E33").Value = "  -1.56%  "

# Row 34 - Monero
$ws.Range("D34").Value = "158.61"
$ws.Range("E34").Value = "  +0.64%  "

# Row 35 - WhiteBITCoin
$ws.Range("D35").Value = "19.04"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36 - USDe
$ws.Range("E36").Value = "  -0.04%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "17.93"
$ws.Range("E37").Value = "  +0.44%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  -4.84%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("D39").Value = "0.297"
$ws.Range("E39").Value = "  -1.86%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "4.32"
$ws.Range("E40").Value = "  -3.12%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "1.48"
$ws.Range("E41").Value = "  -0.06%  "

# Row 42 - ImmutableX
$ws.Range("D42").Value = "1.08"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43 - Aave
$ws.Range("D43").Value = "133.34"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 - was Filecoin, becomes dogwifhat
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").Value = "  -4.50%  "

# Row 45 - was dogwifhat, becomes Filecoin
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  -0.85%  "

# Row 46 - Cronos
$ws.Range("E46").Value = "  -0.80%  "

# Row 47 - ARBITRUM
$ws.Range("D47").Value = "0.480"
$ws.Range("E47").Value = "  -1.33%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "0.557"
$ws.Range("E48").Value = "  -1.29%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.35%  "

# Row 50 - BitgetToken
$ws.Range("D50").Value = "1.12"
$ws.Range("E50").Value = "  +0.04%  "

# Row 51 - Optimism
$ws.Range("D51").Value = "1.38"
$ws.Range("E51").Value = "  -2.57%  "

# Restore the original (default) style now that every value is safely
# stored as text, so no extraneous number-format/style survives on these
# cells.
$dataRange.Style = "Normal"
